# Apply "simulator full-month coverage, persist logs, fix employees" changes.
$wb = $excel.ActiveWorkbook

$wsWeekly = $wb.Worksheets.Item("Weekly Timesheet")
$wsSchema = $wb.Worksheets.Item("Jason Schema")

# --- Fix client names (shared between both sheets: Weekly Timesheet col B, Jason Schema col D) ---
$wsWeekly.Range("B2").Value = "Schauer"
$wsWeekly.Range("B3").Value = "Muncey"
$wsWeekly.Range("B4").Value = "Moulton"
$wsWeekly.Range("B5").Value = "Regan"
$wsWeekly.Range("B6").Value = "Hendricks"

$wsSchema.Range("D2").Value = "Schauer"
$wsSchema.Range("D3").Value = "Muncey"
$wsSchema.Range("D4").Value = "Moulton"
$wsSchema.Range("D5").Value = "Regan"
$wsSchema.Range("D6").Value = "Hendricks"

# --- Simulator full-month coverage: populate Rate (E) and Total (F) for each day row ---
for ($r = 2; $r -le 6; $r++) {
    $wsWeekly.Cells.Item($r, 5).Value = 110
    $wsWeekly.Cells.Item($r, 6).Value = 880
}

# --- Roll the per-day totals up into the subtotal / hourly subtotal / grand total rows ---
$wsWeekly.Range("F8").Value = 4400
$wsWeekly.Range("F11").Value = 4400
$wsWeekly.Range("F13").Value = 4400

# --- Mirror the same Rate/Total figures on the Jason Schema (persisted log) sheet ---
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 6).Value = 110
    $wsSchema.Cells.Item($r, 7).Value = 880
}

# --- Fix the employee id (logged on every row of the Jason Schema sheet) ---
for ($r = 2; $r -le 6; $r++) {
    $wsSchema.Cells.Item($r, 2).Value = "emp_ga4rqytu"
}
